$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "removal rate" (R) inputs that changed ---
$ws.Range("R2").Value = 0.2
$ws.Range("R3").Value = 0.2
$ws.Range("R19").Value = 0.09
$ws.Range("R26").Value = 0.09

# --- Update T and V formulas (survival / count "error" columns) ---
# Individual (non-contiguous / isolated) rows: 2, 3, 12, 16, 18
$ws.Range("T2").Formula  = "=(AVERAGE(J2:K2)/(Q2-R2))-S2"
$ws.Range("V2").Formula  = "=(AVERAGE(H2:I2)/(Q2-R2))-U2"

$ws.Range("T3").Formula  = "=(AVERAGE(J3:K3)/(Q3-R3))-S3"
$ws.Range("V3").Formula  = "=(AVERAGE(H3:I3)/(Q3-R3))-U3"

$ws.Range("T12").Formula = "=(AVERAGE(J12:K12)/(Q12-R12))-S12"
$ws.Range("V12").Formula = "=(AVERAGE(H12:I12)/(Q12-R12))-U12"

$ws.Range("T16").Formula = "=(AVERAGE(J16:K16)/(Q16-R16))-S16"
$ws.Range("V16").Formula = "=(AVERAGE(H16:I16)/(Q16-R16))-U16"

$ws.Range("T18").Formula = "=(AVERAGE(J18:K18)/(Q18-R18))-S18"
$ws.Range("V18").Formula = "=(AVERAGE(H18:I18)/(Q18-R18))-U18"

# Rows 19-21: T set individually, V set together (19:28) so they share one
# contiguous shared-formula group just like the rest of the block below.
$ws.Range("T19").Formula = "=(AVERAGE(J19:K19)/(Q19-R19))-S19"
$ws.Range("T20").Formula = "=(AVERAGE(J20:K20)/(Q20-R20))-S20"
$ws.Range("T21").Formula = "=(AVERAGE(J21:K21)/(Q21-R21))-S21"

# Rows 22-28: same T formula pattern, applied as one contiguous range
$ws.Range("T22:T28").Formula = "=(AVERAGE(J22:K22)/(Q22-R22))-S22"

# V19:V28 all share the identical pattern - apply as a single range so the
# whole block becomes one shared formula group
$ws.Range("V19:V28").Formula = "=(AVERAGE(H19:I19)/(Q19-R19))-U19"

# --- Update the view state: scroll/selection ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("V34").Select()

Write-Host "Edit complete"
